$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: shift the oldest three quarters (columns I/J/K) into the two newly
# added columns (L/M) plus column K, carrying formatting+values in one Copy.
$ws.Range("K7:K102").Copy($ws.Range("M7:M102")) | Out-Null
$ws.Range("J7:J102").Copy($ws.Range("L7:L102")) | Out-Null
$ws.Range("I7:I102").Copy($ws.Range("K7:K102")) | Out-Null

# Step 2: write the refreshed data (two new quarters + restated history) into D:J
$ws.Cells.Item(7,4).Value2 = 43465
$ws.Cells.Item(7,5).Value2 = 43373
$ws.Cells.Item(7,6).Value2 = 43281
$ws.Cells.Item(7,7).Value2 = 43190
$ws.Cells.Item(7,8).Value2 = 43100
$ws.Cells.Item(7,9).Value2 = 43008
$ws.Cells.Item(7,10).Value2 = 42916
$ws.Cells.Item(8,4).Value2 = 1373600
$ws.Cells.Item(8,5).Value2 = 1270100
$ws.Cells.Item(8,6).Value2 = 1261800
$ws.Cells.Item(8,7).Value2 = 1259000
$ws.Cells.Item(8,8).Value2 = 1340700
$ws.Cells.Item(8,9).Value2 = 1315500
$ws.Cells.Item(8,10).Value2 = 1359800
$ws.Cells.Item(9,4).Value2 = "NA"
$ws.Cells.Item(9,5).Value2 = "NA"
$ws.Cells.Item(9,6).Value2 = "NA"
$ws.Cells.Item(9,7).Value2 = "NA"
$ws.Cells.Item(9,8).Value2 = "NA"
$ws.Cells.Item(9,9).Value2 = "NA"
$ws.Cells.Item(9,10).Value2 = "NA"
$ws.Cells.Item(10,4).Value2 = "NA"
$ws.Cells.Item(10,5).Value2 = "NA"
$ws.Cells.Item(10,6).Value2 = "NA"
$ws.Cells.Item(10,7).Value2 = "NA"
$ws.Cells.Item(10,8).Value2 = "NA"
$ws.Cells.Item(10,9).Value2 = "NA"
$ws.Cells.Item(10,10).Value2 = "NA"
$ws.Cells.Item(12,4).Value2 = "NA"
$ws.Cells.Item(12,5).Value2 = "NA"
$ws.Cells.Item(12,6).Value2 = "NA"
$ws.Cells.Item(12,7).Value2 = "NA"
$ws.Cells.Item(12,8).Value2 = "NA"
$ws.Cells.Item(12,9).Value2 = "NA"
$ws.Cells.Item(12,10).Value2 = "NA"
$ws.Cells.Item(13,4).Value2 = 0
$ws.Cells.Item(13,5).Value2 = 0
$ws.Cells.Item(13,6).Value2 = 0
$ws.Cells.Item(13,7).Value2 = 0
$ws.Cells.Item(13,8).Value2 = 0
$ws.Cells.Item(13,9).Value2 = 0
$ws.Cells.Item(13,10).Value2 = 0
$ws.Cells.Item(14,4).Value2 = 0
$ws.Cells.Item(14,5).Value2 = 0
$ws.Cells.Item(14,6).Value2 = 0
$ws.Cells.Item(14,7).Value2 = 0
$ws.Cells.Item(14,8).Value2 = 0
$ws.Cells.Item(14,9).Value2 = 0
$ws.Cells.Item(14,10).Value2 = 0
$ws.Cells.Item(15,4).Value2 = -41800
$ws.Cells.Item(15,5).Value2 = -42200
$ws.Cells.Item(15,6).Value2 = -36100
$ws.Cells.Item(15,7).Value2 = -37900
$ws.Cells.Item(15,8).Value2 = -39300
$ws.Cells.Item(15,9).Value2 = -38400
$ws.Cells.Item(15,10).Value2 = -37300
$ws.Cells.Item(17,4).Value2 = 774000
$ws.Cells.Item(17,5).Value2 = 773300
$ws.Cells.Item(17,6).Value2 = 759000
$ws.Cells.Item(17,7).Value2 = 738900
$ws.Cells.Item(17,8).Value2 = 779400
$ws.Cells.Item(17,9).Value2 = 806100
$ws.Cells.Item(17,10).Value2 = 762700
$ws.Cells.Item(18,4).Value2 = 599600
$ws.Cells.Item(18,5).Value2 = 496800
$ws.Cells.Item(18,6).Value2 = 502800
$ws.Cells.Item(18,7).Value2 = 520100
$ws.Cells.Item(18,8).Value2 = 561300
$ws.Cells.Item(18,9).Value2 = 509400
$ws.Cells.Item(18,10).Value2 = 597000
$ws.Cells.Item(20,4).Value2 = -233800
$ws.Cells.Item(20,5).Value2 = -241100
$ws.Cells.Item(20,6).Value2 = -241000
$ws.Cells.Item(20,7).Value2 = -246400
$ws.Cells.Item(20,8).Value2 = -140700
$ws.Cells.Item(20,9).Value2 = -281100
$ws.Cells.Item(20,10).Value2 = -287800
$ws.Cells.Item(21,4).Value2 = "NA"
$ws.Cells.Item(21,5).Value2 = "NA"
$ws.Cells.Item(21,6).Value2 = "NA"
$ws.Cells.Item(21,7).Value2 = "NA"
$ws.Cells.Item(21,8).Value2 = "NA"
$ws.Cells.Item(21,9).Value2 = "NA"
$ws.Cells.Item(21,10).Value2 = "NA"
$ws.Cells.Item(22,4).Value2 = 0
$ws.Cells.Item(22,5).Value2 = 0
$ws.Cells.Item(22,6).Value2 = 0
$ws.Cells.Item(22,7).Value2 = 0
$ws.Cells.Item(22,8).Value2 = 0
$ws.Cells.Item(22,9).Value2 = 0
$ws.Cells.Item(22,10).Value2 = 0
$ws.Cells.Item(23,4).Value2 = 365800
$ws.Cells.Item(23,5).Value2 = 255700
$ws.Cells.Item(23,6).Value2 = 261900
$ws.Cells.Item(23,7).Value2 = 273700
$ws.Cells.Item(23,8).Value2 = 420600
$ws.Cells.Item(23,9).Value2 = 228300
$ws.Cells.Item(23,10).Value2 = 309300
$ws.Cells.Item(24,4).Value2 = 35200
$ws.Cells.Item(24,5).Value2 = 68300
$ws.Cells.Item(24,6).Value2 = 62400
$ws.Cells.Item(24,7).Value2 = 99600
$ws.Cells.Item(24,8).Value2 = 110600
$ws.Cells.Item(24,9).Value2 = 78500
$ws.Cells.Item(24,10).Value2 = 89900
$ws.Cells.Item(25,4).Value2 = 0
$ws.Cells.Item(25,5).Value2 = 0
$ws.Cells.Item(25,6).Value2 = 0
$ws.Cells.Item(25,7).Value2 = 0
$ws.Cells.Item(25,8).Value2 = 0
$ws.Cells.Item(25,9).Value2 = 0
$ws.Cells.Item(25,10).Value2 = 0
$ws.Cells.Item(26,4).Value2 = 330700
$ws.Cells.Item(26,5).Value2 = 187400
$ws.Cells.Item(26,6).Value2 = 199500
$ws.Cells.Item(26,7).Value2 = 174100
$ws.Cells.Item(26,8).Value2 = 310000
$ws.Cells.Item(26,9).Value2 = 149800
$ws.Cells.Item(26,10).Value2 = 219300
$ws.Cells.Item(27,4).Value2 = 320700
$ws.Cells.Item(27,5).Value2 = 173800
$ws.Cells.Item(27,6).Value2 = 189300
$ws.Cells.Item(27,7).Value2 = 167000
$ws.Cells.Item(27,8).Value2 = 288600
$ws.Cells.Item(27,9).Value2 = 144300
$ws.Cells.Item(27,10).Value2 = 209100
$ws.Cells.Item(28,4).Value2 = 0
$ws.Cells.Item(28,5).Value2 = 0
$ws.Cells.Item(28,6).Value2 = 0
$ws.Cells.Item(28,7).Value2 = 0
$ws.Cells.Item(28,8).Value2 = 0
$ws.Cells.Item(28,9).Value2 = 0
$ws.Cells.Item(28,10).Value2 = 0
$ws.Cells.Item(29,4).Value2 = "NA"
$ws.Cells.Item(29,5).Value2 = "NA"
$ws.Cells.Item(29,6).Value2 = "NA"
$ws.Cells.Item(29,7).Value2 = "NA"
$ws.Cells.Item(29,8).Value2 = "NA"
$ws.Cells.Item(29,9).Value2 = "NA"
$ws.Cells.Item(29,10).Value2 = "NA"
$ws.Cells.Item(30,4).Value2 = 0
$ws.Cells.Item(30,5).Value2 = 0
$ws.Cells.Item(30,6).Value2 = 0
$ws.Cells.Item(30,7).Value2 = 0
$ws.Cells.Item(30,8).Value2 = 0
$ws.Cells.Item(30,9).Value2 = 0
$ws.Cells.Item(30,10).Value2 = 0
$ws.Cells.Item(31,4).Value2 = 0
$ws.Cells.Item(31,5).Value2 = 0
$ws.Cells.Item(31,6).Value2 = 0
$ws.Cells.Item(31,7).Value2 = 0
$ws.Cells.Item(31,8).Value2 = 0
$ws.Cells.Item(31,9).Value2 = 0
$ws.Cells.Item(31,10).Value2 = 0
$ws.Cells.Item(32,4).Value2 = 233800
$ws.Cells.Item(32,5).Value2 = 241100
$ws.Cells.Item(32,6).Value2 = 241000
$ws.Cells.Item(32,7).Value2 = 246400
$ws.Cells.Item(32,8).Value2 = 140700
$ws.Cells.Item(32,9).Value2 = 281100
$ws.Cells.Item(32,10).Value2 = 287800
$ws.Cells.Item(33,4).Value2 = 320700
$ws.Cells.Item(33,5).Value2 = 173800
$ws.Cells.Item(33,6).Value2 = 189300
$ws.Cells.Item(33,7).Value2 = 167000
$ws.Cells.Item(33,8).Value2 = 288600
$ws.Cells.Item(33,9).Value2 = 144300
$ws.Cells.Item(33,10).Value2 = 209100
$ws.Cells.Item(34,4).Value2 = 0
$ws.Cells.Item(34,5).Value2 = 0
$ws.Cells.Item(34,6).Value2 = 0
$ws.Cells.Item(34,7).Value2 = 0
$ws.Cells.Item(34,8).Value2 = 0
$ws.Cells.Item(34,9).Value2 = 0
$ws.Cells.Item(34,10).Value2 = 0
$ws.Cells.Item(35,4).Value2 = 320700
$ws.Cells.Item(35,5).Value2 = 173800
$ws.Cells.Item(35,6).Value2 = 189300
$ws.Cells.Item(35,7).Value2 = 167000
$ws.Cells.Item(35,8).Value2 = 288600
$ws.Cells.Item(35,9).Value2 = 144300
$ws.Cells.Item(35,10).Value2 = 209100
$ws.Cells.Item(38,4).Value2 = 43465
$ws.Cells.Item(38,5).Value2 = 43373
$ws.Cells.Item(38,6).Value2 = 43281
$ws.Cells.Item(38,7).Value2 = 43190
$ws.Cells.Item(38,8).Value2 = 43100
$ws.Cells.Item(38,9).Value2 = 43008
$ws.Cells.Item(38,10).Value2 = 42916
$ws.Cells.Item(41,4).Value2 = 5066600
$ws.Cells.Item(41,5).Value2 = 4435400
$ws.Cells.Item(41,6).Value2 = 4178100
$ws.Cells.Item(41,7).Value2 = 4228700
$ws.Cells.Item(41,8).Value2 = 4967400
$ws.Cells.Item(41,9).Value2 = 4782500
$ws.Cells.Item(41,10).Value2 = 4726900
$ws.Cells.Item(42,4).Value2 = 1517300
$ws.Cells.Item(42,5).Value2 = 1491200
$ws.Cells.Item(42,6).Value2 = 1417500
$ws.Cells.Item(42,7).Value2 = 1607300
$ws.Cells.Item(42,8).Value2 = 1208600
$ws.Cells.Item(42,9).Value2 = 1579900
$ws.Cells.Item(42,10).Value2 = 1911200
$ws.Cells.Item(43,4).Value2 = 0
$ws.Cells.Item(43,5).Value2 = 0
$ws.Cells.Item(43,6).Value2 = 0
$ws.Cells.Item(43,7).Value2 = 0
$ws.Cells.Item(43,8).Value2 = 0
$ws.Cells.Item(43,9).Value2 = 0
$ws.Cells.Item(43,10).Value2 = 0
$ws.Cells.Item(44,4).Value2 = 0
$ws.Cells.Item(44,5).Value2 = 0
$ws.Cells.Item(44,6).Value2 = 0
$ws.Cells.Item(44,7).Value2 = 0
$ws.Cells.Item(44,8).Value2 = 0
$ws.Cells.Item(44,9).Value2 = 0
$ws.Cells.Item(44,10).Value2 = 0
$ws.Cells.Item(45,4).Value2 = 0
$ws.Cells.Item(45,5).Value2 = 0
$ws.Cells.Item(45,6).Value2 = 0
$ws.Cells.Item(45,7).Value2 = 0
$ws.Cells.Item(45,8).Value2 = 0
$ws.Cells.Item(45,9).Value2 = 0
$ws.Cells.Item(45,10).Value2 = 0
$ws.Cells.Item(46,4).Value2 = 0
$ws.Cells.Item(46,5).Value2 = 0
$ws.Cells.Item(46,6).Value2 = 0
$ws.Cells.Item(46,7).Value2 = 0
$ws.Cells.Item(46,8).Value2 = 0
$ws.Cells.Item(46,9).Value2 = 0
$ws.Cells.Item(46,10).Value2 = 0
$ws.Cells.Item(47,4).Value2 = 687900
$ws.Cells.Item(47,5).Value2 = 576800
$ws.Cells.Item(47,6).Value2 = 559900
$ws.Cells.Item(47,7).Value2 = 499200
$ws.Cells.Item(47,8).Value2 = 500800
$ws.Cells.Item(47,9).Value2 = 510500
$ws.Cells.Item(47,10).Value2 = 462500
$ws.Cells.Item(48,4).Value2 = 1632500
$ws.Cells.Item(48,5).Value2 = 1559500
$ws.Cells.Item(48,6).Value2 = 1530300
$ws.Cells.Item(48,7).Value2 = 1512400
$ws.Cells.Item(48,8).Value2 = 1531100
$ws.Cells.Item(48,9).Value2 = 1549200
$ws.Cells.Item(48,10).Value2 = 1517500
$ws.Cells.Item(49,4).Value2 = 2304600
$ws.Cells.Item(49,5).Value2 = 2111200
$ws.Cells.Item(49,6).Value2 = 2080400
$ws.Cells.Item(49,7).Value2 = 1975700
$ws.Cells.Item(49,8).Value2 = 2122100
$ws.Cells.Item(49,9).Value2 = 2088300
$ws.Cells.Item(49,10).Value2 = 2169100
$ws.Cells.Item(50,4).Value2 = 0
$ws.Cells.Item(50,5).Value2 = 0
$ws.Cells.Item(50,6).Value2 = 0
$ws.Cells.Item(50,7).Value2 = 0
$ws.Cells.Item(50,8).Value2 = 0
$ws.Cells.Item(50,9).Value2 = 0
$ws.Cells.Item(50,10).Value2 = 0
$ws.Cells.Item(51,4).Value2 = 0
$ws.Cells.Item(51,5).Value2 = 0
$ws.Cells.Item(51,6).Value2 = 0
$ws.Cells.Item(51,7).Value2 = 0
$ws.Cells.Item(51,8).Value2 = 0
$ws.Cells.Item(51,9).Value2 = 0
$ws.Cells.Item(51,10).Value2 = 0
$ws.Cells.Item(52,4).Value2 = 86800
$ws.Cells.Item(52,5).Value2 = 234300
$ws.Cells.Item(52,6).Value2 = 363500
$ws.Cells.Item(52,7).Value2 = 211200
$ws.Cells.Item(52,8).Value2 = 47600
$ws.Cells.Item(52,9).Value2 = 227700
$ws.Cells.Item(52,10).Value2 = 231600
$ws.Cells.Item(53,4).Value2 = 0
$ws.Cells.Item(53,5).Value2 = 0
$ws.Cells.Item(53,6).Value2 = 0
$ws.Cells.Item(53,7).Value2 = 0
$ws.Cells.Item(53,8).Value2 = 0
$ws.Cells.Item(53,9).Value2 = 0
$ws.Cells.Item(53,10).Value2 = 0
$ws.Cells.Item(54,4).Value2 = 70436400
$ws.Cells.Item(54,5).Value2 = 66129700
$ws.Cells.Item(54,6).Value2 = 65457100
$ws.Cells.Item(54,7).Value2 = 64302600
$ws.Cells.Item(54,8).Value2 = 65250600
$ws.Cells.Item(54,9).Value2 = 65309700
$ws.Cells.Item(54,10).Value2 = 65185400
$ws.Cells.Item(57,4).Value2 = 0
$ws.Cells.Item(57,5).Value2 = 0
$ws.Cells.Item(57,6).Value2 = 0
$ws.Cells.Item(57,7).Value2 = 0
$ws.Cells.Item(57,8).Value2 = 0
$ws.Cells.Item(57,9).Value2 = 0
$ws.Cells.Item(57,10).Value2 = 0
$ws.Cells.Item(58,4).Value2 = 0
$ws.Cells.Item(58,5).Value2 = 0
$ws.Cells.Item(58,6).Value2 = 0
$ws.Cells.Item(58,7).Value2 = 0
$ws.Cells.Item(58,8).Value2 = 0
$ws.Cells.Item(58,9).Value2 = 0
$ws.Cells.Item(58,10).Value2 = 0
$ws.Cells.Item(59,4).Value2 = 53300
$ws.Cells.Item(59,5).Value2 = 182700
$ws.Cells.Item(59,6).Value2 = 129500
$ws.Cells.Item(59,7).Value2 = 131500
$ws.Cells.Item(59,8).Value2 = 51800
$ws.Cells.Item(59,9).Value2 = 274300
$ws.Cells.Item(59,10).Value2 = 208300
$ws.Cells.Item(60,4).Value2 = 0
$ws.Cells.Item(60,5).Value2 = 0
$ws.Cells.Item(60,6).Value2 = 0
$ws.Cells.Item(60,7).Value2 = 0
$ws.Cells.Item(60,8).Value2 = 0
$ws.Cells.Item(60,9).Value2 = 0
$ws.Cells.Item(60,10).Value2 = 0
$ws.Cells.Item(61,4).Value2 = 11906900
$ws.Cells.Item(61,5).Value2 = 11129400
$ws.Cells.Item(61,6).Value2 = 10589700
$ws.Cells.Item(61,7).Value2 = 10071700
$ws.Cells.Item(61,8).Value2 = 10897200
$ws.Cells.Item(61,9).Value2 = 11770200
$ws.Cells.Item(61,10).Value2 = 11640800
$ws.Cells.Item(62,4).Value2 = 704400
$ws.Cells.Item(62,5).Value2 = 690900
$ws.Cells.Item(62,6).Value2 = 844100
$ws.Cells.Item(62,7).Value2 = 658100
$ws.Cells.Item(62,8).Value2 = 717000
$ws.Cells.Item(62,9).Value2 = 644200
$ws.Cells.Item(62,10).Value2 = 829300
$ws.Cells.Item(63,4).Value2 = 0
$ws.Cells.Item(63,5).Value2 = 0
$ws.Cells.Item(63,6).Value2 = 0
$ws.Cells.Item(63,7).Value2 = 0
$ws.Cells.Item(63,8).Value2 = 0
$ws.Cells.Item(63,9).Value2 = 0
$ws.Cells.Item(63,10).Value2 = 0
$ws.Cells.Item(64,4).Value2 = 0
$ws.Cells.Item(64,5).Value2 = 0
$ws.Cells.Item(64,6).Value2 = 0
$ws.Cells.Item(64,7).Value2 = 0
$ws.Cells.Item(64,8).Value2 = 0
$ws.Cells.Item(64,9).Value2 = 0
$ws.Cells.Item(64,10).Value2 = 0
$ws.Cells.Item(65,4).Value2 = 0
$ws.Cells.Item(65,5).Value2 = 0
$ws.Cells.Item(65,6).Value2 = 0
$ws.Cells.Item(65,7).Value2 = 0
$ws.Cells.Item(65,8).Value2 = 0
$ws.Cells.Item(65,9).Value2 = 0
$ws.Cells.Item(65,10).Value2 = 0
$ws.Cells.Item(66,4).Value2 = 62484700
$ws.Cells.Item(66,5).Value2 = 58721300
$ws.Cells.Item(66,6).Value2 = 58180500
$ws.Cells.Item(66,7).Value2 = 57277900
$ws.Cells.Item(66,8).Value2 = 57854500
$ws.Cells.Item(66,9).Value2 = 58232300
$ws.Cells.Item(66,10).Value2 = 58199500
$ws.Cells.Item(68,4).Value2 = 0
$ws.Cells.Item(68,5).Value2 = 0
$ws.Cells.Item(68,6).Value2 = 0
$ws.Cells.Item(68,7).Value2 = 0
$ws.Cells.Item(68,8).Value2 = 0
$ws.Cells.Item(68,9).Value2 = 0
$ws.Cells.Item(68,10).Value2 = 0
$ws.Cells.Item(69,4).Value2 = 0
$ws.Cells.Item(69,5).Value2 = 0
$ws.Cells.Item(69,6).Value2 = 0
$ws.Cells.Item(69,7).Value2 = 0
$ws.Cells.Item(69,8).Value2 = 0
$ws.Cells.Item(69,9).Value2 = 0
$ws.Cells.Item(69,10).Value2 = 0
$ws.Cells.Item(70,4).Value2 = 0
$ws.Cells.Item(70,5).Value2 = 0
$ws.Cells.Item(70,6).Value2 = 0
$ws.Cells.Item(70,7).Value2 = 0
$ws.Cells.Item(70,8).Value2 = 0
$ws.Cells.Item(70,9).Value2 = 0
$ws.Cells.Item(70,10).Value2 = 0
$ws.Cells.Item(71,4).Value2 = 0
$ws.Cells.Item(71,5).Value2 = 0
$ws.Cells.Item(71,6).Value2 = 0
$ws.Cells.Item(71,7).Value2 = 0
$ws.Cells.Item(71,8).Value2 = 0
$ws.Cells.Item(71,9).Value2 = 0
$ws.Cells.Item(71,10).Value2 = 0
$ws.Cells.Item(72,4).Value2 = 5218400
$ws.Cells.Item(72,5).Value2 = 4873100
$ws.Cells.Item(72,6).Value2 = 4765600
$ws.Cells.Item(72,7).Value2 = 4594600
$ws.Cells.Item(72,8).Value2 = 4873100
$ws.Cells.Item(72,9).Value2 = 4572200
$ws.Cells.Item(72,10).Value2 = 4427900
$ws.Cells.Item(73,4).Value2 = 0
$ws.Cells.Item(73,5).Value2 = 0
$ws.Cells.Item(73,6).Value2 = 0
$ws.Cells.Item(73,7).Value2 = 0
$ws.Cells.Item(73,8).Value2 = 0
$ws.Cells.Item(73,9).Value2 = 0
$ws.Cells.Item(73,10).Value2 = 0
$ws.Cells.Item(74,4).Value2 = 0
$ws.Cells.Item(74,5).Value2 = 0
$ws.Cells.Item(74,6).Value2 = 0
$ws.Cells.Item(74,7).Value2 = 0
$ws.Cells.Item(74,8).Value2 = 0
$ws.Cells.Item(74,9).Value2 = 0
$ws.Cells.Item(74,10).Value2 = 0
$ws.Cells.Item(75,4).Value2 = 0
$ws.Cells.Item(75,5).Value2 = 0
$ws.Cells.Item(75,6).Value2 = 0
$ws.Cells.Item(75,7).Value2 = 0
$ws.Cells.Item(75,8).Value2 = 0
$ws.Cells.Item(75,9).Value2 = 0
$ws.Cells.Item(75,10).Value2 = 0
$ws.Cells.Item(76,4).Value2 = 7951700
$ws.Cells.Item(76,5).Value2 = 7408400
$ws.Cells.Item(76,6).Value2 = 7276600
$ws.Cells.Item(76,7).Value2 = 7024700
$ws.Cells.Item(76,8).Value2 = 7396100
$ws.Cells.Item(76,9).Value2 = 7077400
$ws.Cells.Item(76,10).Value2 = 6985900
$ws.Cells.Item(77,4).Value2 = 0
$ws.Cells.Item(77,5).Value2 = 0
$ws.Cells.Item(77,6).Value2 = 0
$ws.Cells.Item(77,7).Value2 = 0
$ws.Cells.Item(77,8).Value2 = 0
$ws.Cells.Item(77,9).Value2 = 0
$ws.Cells.Item(77,10).Value2 = 0
$ws.Cells.Item(80,4).Value2 = 43465
$ws.Cells.Item(80,5).Value2 = 43373
$ws.Cells.Item(80,6).Value2 = 43281
$ws.Cells.Item(80,7).Value2 = 43190
$ws.Cells.Item(80,8).Value2 = 43100
$ws.Cells.Item(80,9).Value2 = 43008
$ws.Cells.Item(80,10).Value2 = 42916
$ws.Cells.Item(81,4).Value2 = 320700
$ws.Cells.Item(81,5).Value2 = 173800
$ws.Cells.Item(81,6).Value2 = 189300
$ws.Cells.Item(81,7).Value2 = 167000
$ws.Cells.Item(81,8).Value2 = 288600
$ws.Cells.Item(81,9).Value2 = 144300
$ws.Cells.Item(81,10).Value2 = 209100
$ws.Cells.Item(83,4).Value2 = 0
$ws.Cells.Item(83,5).Value2 = 0
$ws.Cells.Item(83,6).Value2 = 0
$ws.Cells.Item(83,7).Value2 = 0
$ws.Cells.Item(83,8).Value2 = 0
$ws.Cells.Item(83,9).Value2 = 0
$ws.Cells.Item(83,10).Value2 = 0
$ws.Cells.Item(84,4).Value2 = 0
$ws.Cells.Item(84,5).Value2 = 0
$ws.Cells.Item(84,6).Value2 = 0
$ws.Cells.Item(84,7).Value2 = 0
$ws.Cells.Item(84,8).Value2 = 0
$ws.Cells.Item(84,9).Value2 = 0
$ws.Cells.Item(84,10).Value2 = 0
$ws.Cells.Item(85,4).Value2 = 0
$ws.Cells.Item(85,5).Value2 = 0
$ws.Cells.Item(85,6).Value2 = 0
$ws.Cells.Item(85,7).Value2 = 0
$ws.Cells.Item(85,8).Value2 = 0
$ws.Cells.Item(85,9).Value2 = 0
$ws.Cells.Item(85,10).Value2 = 0
$ws.Cells.Item(86,4).Value2 = 0
$ws.Cells.Item(86,5).Value2 = 0
$ws.Cells.Item(86,6).Value2 = 0
$ws.Cells.Item(86,7).Value2 = 0
$ws.Cells.Item(86,8).Value2 = 0
$ws.Cells.Item(86,9).Value2 = 0
$ws.Cells.Item(86,10).Value2 = 0
$ws.Cells.Item(87,4).Value2 = 0
$ws.Cells.Item(87,5).Value2 = 0
$ws.Cells.Item(87,6).Value2 = 0
$ws.Cells.Item(87,7).Value2 = 0
$ws.Cells.Item(87,8).Value2 = 0
$ws.Cells.Item(87,9).Value2 = 0
$ws.Cells.Item(87,10).Value2 = 0
$ws.Cells.Item(88,4).Value2 = 0
$ws.Cells.Item(88,5).Value2 = 0
$ws.Cells.Item(88,6).Value2 = 0
$ws.Cells.Item(88,7).Value2 = 0
$ws.Cells.Item(88,8).Value2 = 0
$ws.Cells.Item(88,9).Value2 = 0
$ws.Cells.Item(88,10).Value2 = 0
$ws.Cells.Item(89,4).Value2 = 0
$ws.Cells.Item(89,5).Value2 = 0
$ws.Cells.Item(89,6).Value2 = 0
$ws.Cells.Item(89,7).Value2 = 0
$ws.Cells.Item(89,8).Value2 = 0
$ws.Cells.Item(89,9).Value2 = 0
$ws.Cells.Item(89,10).Value2 = 0
$ws.Cells.Item(91,4).Value2 = 0
$ws.Cells.Item(91,5).Value2 = 0
$ws.Cells.Item(91,6).Value2 = 0
$ws.Cells.Item(91,7).Value2 = 0
$ws.Cells.Item(91,8).Value2 = 0
$ws.Cells.Item(91,9).Value2 = 0
$ws.Cells.Item(91,10).Value2 = 0
$ws.Cells.Item(92,4).Value2 = 0
$ws.Cells.Item(92,5).Value2 = 0
$ws.Cells.Item(92,6).Value2 = 0
$ws.Cells.Item(92,7).Value2 = 0
$ws.Cells.Item(92,8).Value2 = 0
$ws.Cells.Item(92,9).Value2 = 0
$ws.Cells.Item(92,10).Value2 = 0
$ws.Cells.Item(93,4).Value2 = 0
$ws.Cells.Item(93,5).Value2 = 0
$ws.Cells.Item(93,6).Value2 = 0
$ws.Cells.Item(93,7).Value2 = 0
$ws.Cells.Item(93,8).Value2 = 0
$ws.Cells.Item(93,9).Value2 = 0
$ws.Cells.Item(93,10).Value2 = 0
$ws.Cells.Item(94,4).Value2 = 0
$ws.Cells.Item(94,5).Value2 = 0
$ws.Cells.Item(94,6).Value2 = 0
$ws.Cells.Item(94,7).Value2 = 0
$ws.Cells.Item(94,8).Value2 = 0
$ws.Cells.Item(94,9).Value2 = 0
$ws.Cells.Item(94,10).Value2 = 0
$ws.Cells.Item(96,4).Value2 = 0
$ws.Cells.Item(96,5).Value2 = 0
$ws.Cells.Item(96,6).Value2 = 0
$ws.Cells.Item(96,7).Value2 = 0
$ws.Cells.Item(96,8).Value2 = 0
$ws.Cells.Item(96,9).Value2 = 0
$ws.Cells.Item(96,10).Value2 = 0
$ws.Cells.Item(97,4).Value2 = 0
$ws.Cells.Item(97,5).Value2 = 0
$ws.Cells.Item(97,6).Value2 = 0
$ws.Cells.Item(97,7).Value2 = 0
$ws.Cells.Item(97,8).Value2 = 0
$ws.Cells.Item(97,9).Value2 = 0
$ws.Cells.Item(97,10).Value2 = 0
$ws.Cells.Item(98,4).Value2 = 0
$ws.Cells.Item(98,5).Value2 = 0
$ws.Cells.Item(98,6).Value2 = 0
$ws.Cells.Item(98,7).Value2 = 0
$ws.Cells.Item(98,8).Value2 = 0
$ws.Cells.Item(98,9).Value2 = 0
$ws.Cells.Item(98,10).Value2 = 0
$ws.Cells.Item(99,4).Value2 = 0
$ws.Cells.Item(99,5).Value2 = 0
$ws.Cells.Item(99,6).Value2 = 0
$ws.Cells.Item(99,7).Value2 = 0
$ws.Cells.Item(99,8).Value2 = 0
$ws.Cells.Item(99,9).Value2 = 0
$ws.Cells.Item(99,10).Value2 = 0
$ws.Cells.Item(100,4).Value2 = 0
$ws.Cells.Item(100,5).Value2 = 0
$ws.Cells.Item(100,6).Value2 = 0
$ws.Cells.Item(100,7).Value2 = 0
$ws.Cells.Item(100,8).Value2 = 0
$ws.Cells.Item(100,9).Value2 = 0
$ws.Cells.Item(100,10).Value2 = 0
$ws.Cells.Item(101,4).Value2 = 0
$ws.Cells.Item(101,5).Value2 = 0
$ws.Cells.Item(101,6).Value2 = 0
$ws.Cells.Item(101,7).Value2 = 0
$ws.Cells.Item(101,8).Value2 = 0
$ws.Cells.Item(101,9).Value2 = 0
$ws.Cells.Item(101,10).Value2 = 0
$ws.Cells.Item(102,4).Value2 = 0
$ws.Cells.Item(102,5).Value2 = 0
$ws.Cells.Item(102,6).Value2 = 0
$ws.Cells.Item(102,7).Value2 = 0
$ws.Cells.Item(102,8).Value2 = 0
$ws.Cells.Item(102,9).Value2 = 0
$ws.Cells.Item(102,10).Value2 = 0
